$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.381.33"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.610.54"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.68"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.65"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.609.73"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.28"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.081.61"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000182"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.132.84"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.606.50"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "371.97"
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.04"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.75"
$ws.Range("E23").Value = "  -4.78%  "
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.67"
$ws.Range("E25").Value = "  +4.94%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.742.86"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "580.29"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0986"
$ws.Range("E31").Value = "  -6.63%  "
$ws.Range("E32").Value = "  -5.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.67"
$ws.Range("E33").Value = "  -3.66%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.53"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  -4.06%  "
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.14"
$ws.Range("E44").Value = "  +4.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "153.12"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0283"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0778"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.68"
$ws.Range("E50").Value = "  -6.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.38"
$ws.Range("E51").Value = "  +1.17%  "
